$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.620.60'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '1.577.18'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.45'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.84'
$ws.Range("E8").Value = '  +1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.03'
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("E10").Value = '  -1.11%  '
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '1.803.74'
$ws.Range("D14").Value = '1.565.22'
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").Value = '28.629.00'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.68'
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.27'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '230.75'
$ws.Range("E19").Value = '  +1.02%  '
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.88'
$ws.Range("E23").Value = '  -4.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.17'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.09'
$ws.Range("E25").Value = '  +6.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.54'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.45'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0483'
$ws.Range("E31").Value = '  +2.36%  '
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("D35").Value = '1.399.70'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  +3.99%  '
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("E39").Value = '  +2.91%  '
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.793'
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.48'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '62.99'
$ws.Range("E48").Value = '  -1.92%  '
$ws.Range("D49").Value = '1.715.05'
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.55'
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("D51").Value = '0.0₆0103'
$ws.Range("E51").Value = '  -0.63%  '
